# nameList.xlsx update — "Add files via upload"
#
# Sheet1 holds a name-pairing list. This update:
#   - fixes two garbled/typo'd Chinese names (D3, C9)
#   - shifts several of the paired-name columns (E/F/G) up, adding a
#     newly paired name ("Chris Ng") and clearing now-empty trailing cells
#   - fills in the previously-empty H10:H12 cells
#   - leaves the active selection on E12

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Try to restore the last-saved window size recorded in the workbook view.
# (Harmless no-op if the host doesn't round-trip this cosmetic property.)
try {
    $win = $excel.ActiveWindow
    $win.Width = 18615
    $win.Height = 10830
} catch {}

# --- Corrected names -------------------------------------------------
$ws.Range("C9").Value = "陳小凡"
$ws.Range("D3").Value = "廖素琼"

# --- Column F (shift up by one row, rows 1-11; clear row 12) --------
$ws.Range("F1").Value = "*Lam Kuen"

# F2 also drops its (no-op) explicit fill formatting in the source file.
$ws.Range("F2").ClearFormats()
$ws.Range("F2").Value = "Faye"

$ws.Range("F3").Value = "Rohda"
$ws.Range("F4").Value = "Patrick"
$ws.Range("F5").Value = "Leo"
$ws.Range("F6").Value = "Fruit"
$ws.Range("F7").Value = "Chan Dan"
$ws.Range("F8").Value = "Ho Kim Chin"
$ws.Range("F9").Value = "Ho Ming"
$ws.Range("F10").Value = "Harry Cheung"
$ws.Range("F11").Value = "Harry Cheung Wife"
$ws.Range("F12").ClearContents()

# --- Column H (newly populated rows 10-12) ---------------------------
$ws.Range("H10").Value = "*Anthony "
$ws.Range("H11").Value = "Tim"
$ws.Range("H12").Value = "Patrick"

# --- Column E (rows 11-13 shift; rows 18-21 rearranged) --------------
$ws.Range("E11").Value = "Mak Wing"
$ws.Range("E12").Value = "Chris Ng"
$ws.Range("E13").ClearContents()

$ws.Range("E19").Value = "?Wong Kei"
$ws.Range("E20").Value = "Zuey Tsui"
$ws.Range("E21").Value = "Guanglei"

# --- Column G (rows 18-22 shift up by two; clear trailing rows) ------
$ws.Range("G18").Value = "Bean Man"
$ws.Range("G19").Value = "Foo Kwai"
$ws.Range("G20").Value = "Ellen "
$ws.Range("G21").ClearContents()
$ws.Range("G22").ClearContents()

# --- Final selection ---------------------------------------------------
$ws.Range("E12").Select()
